$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The "Noted:" signature block used the field placeholders ${cscAdviser} and
# ${oicOsa} side by side. The edit repurposes that line to present
# ${cscPresident} first (shifting the old adviser placeholder into the slot
# formerly occupied by ${oicOsa}) and tightens the spacing that separates
# the two placeholders.
# ---------------------------------------------------------------------------

# 1) ${cscAdviser} -> ${cscPresident}
$rngAdviser = $d.Content
$rngAdviser.Find.Execute("cscAdviser", $false, $false, $false, $false, $false, `
    $true, 1, $false, "cscPresident", 2)

# 2) ${oicOsa} -> ${cscAdviser}
$rngOsa = $d.Content
$rngOsa.Find.Execute("oicOsa", $false, $false, $false, $false, $false, `
    $true, 1, $false, "cscAdviser", 2)

# 3) Shrink the run of spaces that pads out to the second placeholder from
#    45 characters down to 13.
$spaces45 = "".PadLeft(45, ' ')
$spaces13 = "".PadLeft(13, ' ')
$rngSpacing = $d.Content
$rngSpacing.Find.Execute($spaces45, $false, $false, $false, $false, $false, `
    $true, 1, $false, $spaces13, 2)

# ---------------------------------------------------------------------------
# The "CSC President ... Adviser" line had a long run of spaces pushing the
# word "Adviser" out under the signature line; collapse it to a single
# space.
# ---------------------------------------------------------------------------
$spaces20 = "".PadLeft(20, ' ')
$rngLabel = $d.Content
$rngLabel.Find.Execute(($spaces20 + "Adviser"), $false, $false, $false, $false, $false, `
    $true, 1, $false, (" Adviser"), 2)
